# Update the sheet with newly re-computed TPM-based values (Lama1-Rpsa L-R pair).
# Only the "Receptor" (Rpsa) derived columns (M-T) are affected, since the
# underlying per-cluster Rpsa TPM values changed while Lama1 (ligand) values
# stayed the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sending: FAPs, Target: ECs)
$ws.Range("M2").Value = 91.74689966666665
$ws.Range("N2").Value = 275.2406989999999
$ws.Range("O2").Value = 0.1908387282982634
$ws.Range("P2").Value = 0.1908387282982634
$ws.Range("Q2").Value = 5.162200474344777
$ws.Range("R2").Value = 46.45980426910299
$ws.Range("S2").Value = 0.05207278326675228
$ws.Range("T2").Value = 0.05207278326675228

# Row 3 (Sending: FAPs, Target: FAPs) - receptor average/total unchanged, specificity shifted
$ws.Range("O3").Value = 0.296899627499751
$ws.Range("P3").Value = 0.296899627499751
$ws.Range("S3").Value = 0.08101285358918787
$ws.Range("T3").Value = 0.08101285358918787

# Row 4 (Sending: FAPs, Target: MuSCs)
$ws.Range("M4").Value = 167.6324513333334
$ws.Range("N4").Value = 502.8973540000001
$ws.Range("O4").Value = 0.348684957750095
$ws.Range("P4").Value = 0.348684957750095
$ws.Range("Q4").Value = 9.431951629237558
$ws.Range("R4").Value = 84.88756466313801
$ws.Range("S4").Value = 0.09514314204043352
$ws.Range("T4").Value = 0.09514314204043352

# Row 5 (Sending: FAPs, Target: Resolving-Mac)
$ws.Range("M5").Value = 78.64050433333334
$ws.Range("N5").Value = 235.921513
$ws.Range("O5").Value = 0.1635766864518907
$ws.Range("P5").Value = 0.1635766864518907
$ws.Range("Q5").Value = 4.424760403317889
$ws.Range("R5").Value = 39.822843629861
$ws.Range("S5").Value = 0.04463398712126248
$ws.Range("T5").Value = 0.04463398712126248

# Row 6 (Sending: MuSCs, Target: ECs)
$ws.Range("M6").Value = 91.74689966666665
$ws.Range("N6").Value = 275.2406989999999
$ws.Range("O6").Value = 0.1908387282982634
$ws.Range("P6").Value = 0.1908387282982634
$ws.Range("Q6").Value = 13.75646897142022
$ws.Range("R6").Value = 123.808220742782
$ws.Range("S6").Value = 0.1387659450315111
$ws.Range("T6").Value = 0.1387659450315111

# Row 7 (Sending: MuSCs, Target: FAPs) - receptor average/total unchanged, specificity shifted
$ws.Range("O7").Value = 0.296899627499751
$ws.Range("P7").Value = 0.296899627499751
$ws.Range("S7").Value = 0.2158867739105631
$ws.Range("T7").Value = 0.2158867739105631

# Row 8 (Sending: MuSCs, Target: MuSCs)
$ws.Range("M8").Value = 167.6324513333334
$ws.Range("N8").Value = 502.8973540000001
$ws.Range("O8").Value = 0.348684957750095
$ws.Range("P8").Value = 0.348684957750095
$ws.Range("Q8").Value = 25.13469799795245
$ws.Range("R8").Value = 226.2122819815721
$ws.Range("S8").Value = 0.2535418157096615
$ws.Range("T8").Value = 0.2535418157096615

# Row 9 (Sending: MuSCs, Target: Resolving-Mac)
$ws.Range("M9").Value = 78.64050433333334
$ws.Range("N9").Value = 235.921513
$ws.Range("O9").Value = 0.1635766864518907
$ws.Range("P9").Value = 0.1635766864518907
$ws.Range("Q9").Value = 11.79130479273711
$ws.Range("R9").Value = 106.121743134634
$ws.Range("S9").Value = 0.1189426993306282
$ws.Range("T9").Value = 0.1189426993306282
